# Updates the crypto price/volume table to the latest scraped snapshot.
# D (Price) and E (Volume 1h) columns hold numeric-looking values that are
# stored as TEXT in this sheet (not numbers/percentages), so each literal is
# written with a leading apostrophe (PowerShell '' escapes it) to force Excel
# to keep it as text instead of auto-converting it to a number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '''305.71'
$ws.Range("E2").Value = '''0.59%'

# Row 3
$ws.Range("D3").Value = '''36.18'
$ws.Range("E3").Value = '''-1.65%'

# Row 4
$ws.Range("D4").Value = '''5.036'

# Row 5
$ws.Range("D5").Value = '''0.07877'
$ws.Range("E5").Value = '''1.40%'

# Row 6
$ws.Range("D6").Value = '''2.257'
$ws.Range("E6").Value = '''7.36%'

# Row 7
$ws.Range("D7").Value = '''7.991'
$ws.Range("E7").Value = '''-0.51%'

# Row 8
$ws.Range("B8").Value = 'GateToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D8").Value = '''4.149'
$ws.Range("E8").Value = '''2.41%'

# Row 9
$ws.Range("B9").Value = 'MXToken'
$ws.Range("C9").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D9").Value = '''0.9277'
$ws.Range("E9").Value = '''0.73%'

# Row 10
$ws.Range("B10").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C10").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D10").Value = '''0.09804'
$ws.Range("E10").Value = '''-0.69%'

# Row 11
$ws.Range("B11").Value = 'WazirX'
$ws.Range("C11").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D11").Value = '''0.1866'
$ws.Range("E11").Value = '''0.22%'

# Row 12
$ws.Range("B12").Value = 'MandalaExchangeToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D12").Value = '''0.08975'
$ws.Range("E12").Value = '''3.70%'

# Row 13
$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D13").Value = '''0.03754'
$ws.Range("E13").Value = '''4.52%'

# Row 14
$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D14").Value = '''0.09920'
$ws.Range("E14").Value = '''-0.57%'

# Row 15
$ws.Range("B15").Value = 'BitForexToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D15").Value = '''0.001446'
$ws.Range("E15").Value = '''-2.32%'

# Row 16
$ws.Range("B16").Value = 'TigerCash'
$ws.Range("C16").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D16").Value = '''0.005704'
$ws.Range("E16").Value = '''0.22%'

# Row 17
$ws.Range("B17").Value = 'LEO'
$ws.Range("C17").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D17").Value = '''3.460'
$ws.Range("E17").Value = '''0.01%'

# Row 18
$ws.Range("E18").Value = '''5.01%'

# Row 19
$ws.Range("E19").Value = '''-2.16%'

# Row 20
$ws.Range("E20").Value = '''-0.98%'

# Row 21
$ws.Range("D21").Value = '''5.111'
$ws.Range("E21").Value = '''3.36%'

# Row 22
$ws.Range("D22").Value = '''0.2252'
$ws.Range("E22").Value = '''1.90%'

# Row 23
$ws.Range("D23").Value = '''0.04576'
$ws.Range("E23").Value = '''-0.63%'

# Row 24
$ws.Range("D24").Value = '''0.001234'
$ws.Range("E24").Value = '''-0.28%'

# Row 25
$ws.Range("D25").Value = '''0.004769'
$ws.Range("E25").Value = '''-7.17%'

# Row 26
$ws.Range("D26").Value = '''0.0001303'
$ws.Range("E26").Value = '''-7.52%'

# Row 39
$ws.Range("D39").Value = '''0.01936'
$ws.Range("E39").Value = '''8.07%'

# Row 40
$ws.Range("D40").Value = '''0.04934'
$ws.Range("E40").Value = '''5.63%'

# Row 41
$ws.Range("D41").Value = '''0.007807'
$ws.Range("E41").Value = '''1.25%'

# Row 42
$ws.Range("E42").Value = '''-0.48%'

# Row 43
$ws.Range("E43").Value = '''2.39%'

# Row 44
$ws.Range("E44").Value = '''-3.58%'

# Row 45
$ws.Range("D45").Value = '''0.01142'
$ws.Range("E45").Value = '''9.55%'

# Row 46
$ws.Range("D46").Value = '''0.00006163'
$ws.Range("E46").Value = '''-2.49%'

# Row 47
$ws.Range("E47").Value = '''-0.35%'

# Row 48
$ws.Range("E48").Value = '''52.91%'

# Row 49
$ws.Range("E49").Value = '''-10.31%'

# Row 50
$ws.Range("D50").Value = '''0.00002103'
$ws.Range("E50").Value = '''-0.35%'

# Row 51
$ws.Range("D51").Value = '''0.0002003'
$ws.Range("E51").Value = '''-0.35%'
